# Auto update md files
#
# Refresh the weekly task tracker: a handful of rows had their "当前状态"
# (current status) updated, and two tasks slipped to new due dates.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 - 牙医检查 (dentist checkup): 进行中 -> 未完成
$ws.Range("E3").Value = "未完成"

# Row 4 - 提交健身房会员申请 (gym membership): 已完成 -> 进行中
$ws.Range("E4").Value = "进行中"

# Row 5 - 准备周五部门总结PPT (Friday summary deck): due date pushed out,
# and status slips back from 进行中 to 未完成
$ws.Range("B5").Value = 46051
$ws.Range("E5").Value = "未完成"

# Row 6 - 超市采购生活用品 (grocery run): due date moved up, and the task
# actually got finished: 未完成 -> 已完成
$ws.Range("B6").Value = 46050
$ws.Range("E6").Value = "已完成"

# Row 7 - 整理书房电子发票 (sorting receipts): 未完成 -> 进行中
$ws.Range("E7").Value = "进行中"

# Leave the cursor on the last cell that was touched, as Excel would.
$ws.Range("E5").Select() | Out-Null
